$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.101.40"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "2.353.49"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.20"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.74"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("E7").Value = "  +0.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.564"
$ws.Range("E8").Value = "  +4.97%  "
$ws.Range("E9").Value = "  +3.00%  "
$ws.Range("E10").Value = "  +2.66%  "
$ws.Range("E11").Value = "  -2.02%  "
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").Value = "2.770.71"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.76"
$ws.Range("D15").Value = "58.057.49"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("D17").Value = "2.348.15"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.82"
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("E19").Value = "  +2.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "329.08"
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.75"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "63.37"
$ws.Range("E23").Value = "  +3.09%  "
$ws.Range("E24").Value = "  -2.58%  "
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.24"
$ws.Range("E26").Value = "  -2.78%  "
$ws.Range("E27").Value = "  -4.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.75"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.24"
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").Value = "0.0₃0736"
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.12"
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.35"
$ws.Range("E32").Value = "  -0.71%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("B34").Value = "SuiNetwork"
$ws.Range("C34").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -2.61%  "
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("E37").Value = "  -2.64%  "
$ws.Range("E38").Value = "  -2.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.381"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "140.98"
$ws.Range("E40").Value = "  -5.65%  "
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "288.90"
$ws.Range("E42").Value = "  +0.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0948"
$ws.Range("E43").Value = "  +2.37%  "
$ws.Range("E44").Value = "  +1.97%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.566"
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.89"
$ws.Range("E46").Value = "  -1.86%  "
$ws.Range("E47").Value = "  +2.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.08"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("E51").Value = "  +0.55%  "

Write-Output "Applied all cell updates"